# Auto-generated edit script: updates crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.850.39"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "2.549.91"
$ws.Range("E3").Value = "  +5.70%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'576.19"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").Value = "'149.29"
$ws.Range("E6").Value = "  +7.81%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "2.549.62"
$ws.Range("E9").Value = "  +5.77%  "
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").Value = "'0.360"
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("E14").Value = "  +9.23%  "
$ws.Range("D15").Value = "3.007.00"
$ws.Range("E15").Value = "  +5.63%  "
$ws.Range("D16").Value = "63.654.36"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("D18").Value = "2.549.41"
$ws.Range("E18").Value = "  +5.23%  "
$ws.Range("D19").Value = "'11.62"
$ws.Range("E19").Value = "  +4.78%  "
$ws.Range("D20").Value = "'345.10"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'4.38"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").Value = "'6.91"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'66.20"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  +3.57%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'8.35"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "'1.43"
$ws.Range("E29").Value = "  +4.37%  "
$ws.Range("D30").Value = "0.0₃0839"
$ws.Range("E30").Value = "  +7.38%  "
$ws.Range("D31").Value = "'1.89"
$ws.Range("E31").Value = "  +4.81%  "
$ws.Range("D32").Value = "'6.88"
$ws.Range("E32").Value = "  +8.11%  "
$ws.Range("D33").Value = "'176.75"
$ws.Range("E33").Value = "  +3.70%  "
$ws.Range("D34").Value = "'1.60"
$ws.Range("E34").Value = "  +13.68%  "
$ws.Range("D35").Value = "'424.23"
$ws.Range("E35").Value = "  +12.94%  "
$ws.Range("D36").Value = "'0.406"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D37").Value = "'19.23"
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("D43").Value = "'154.34"
$ws.Range("E43").Value = "  +7.01%  "
$ws.Range("D44").Value = "'3.81"
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("D45").Value = "'21.08"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").Value = "'0.614"
$ws.Range("E46").Value = "  +4.90%  "
$ws.Range("D47").Value = "'0.0535"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("D48").Value = "'0.0970"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("D49").Value = "'19.12"
$ws.Range("E49").Value = "  +6.16%  "
$ws.Range("D50").Value = "'0.0233"
$ws.Range("E50").Value = "  +5.71%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0233"
$ws.Range("E51").Value = "  +8.71%  "
